# Applies the leve-profit recalculation updates captured in the commit diff
# (currentAveragePrice(NQ/HQ) / LevePrice(NQ/HQ) / LeveProfit(NQ/HQ), cols H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1645.4546
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1766.6666
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1766.6666
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2116.6666
$ws.Range("H86").Value = 500951.5
$ws.Range("I86").Value = 500951.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 500951.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -499828.5
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 500951.5
$ws.Range("I89").Value = 500951.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 2504757.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2499141.5
$ws.Range("N89").ClearContents()
$ws.Range("H116").Value = 2151.8667
$ws.Range("I116").Value = 1296.6666
$ws.Range("J116").Value = 3007.0667
$ws.Range("K116").Value = 1296.6666
$ws.Range("L116").Value = 3007.0667
$ws.Range("M116").Value = 2145.3334
$ws.Range("N116").Value = -9891.066699999999
$ws.Range("H132").Value = 6805729.5
$ws.Range("I132").Value = 10207767
$ws.Range("J132").Value = 1654.2858
$ws.Range("K132").Value = 30623301
$ws.Range("L132").Value = 4962.857400000001
$ws.Range("M132").Value = -30620771
$ws.Range("N132").Value = -10022.8574
$ws.Range("H137").Value = 1615.2174
$ws.Range("I137").Value = 1271.5625
$ws.Range("K137").Value = 3814.6875
$ws.Range("M137").Value = -1264.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27778452
$ws.Range("I2").Value = 52632016
$ws.Range("J2").Value = 941.35297
$ws.Range("K2").Value = 52632016
$ws.Range("L2").Value = 941.35297
$ws.Range("M2").Value = -52631903
$ws.Range("N2").Value = -1167.35297
$ws.Range("H32").Value = 21539.54
$ws.Range("I32").Value = 22814.809
$ws.Range("J32").Value = 16438.46
$ws.Range("K32").Value = 22814.809
$ws.Range("L32").Value = 16438.46
$ws.Range("M32").Value = -22527.809
$ws.Range("N32").Value = -17012.46
$ws.Range("H45").Value = 23810854
$ws.Range("I45").Value = 37038270
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 37038270
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -37037893
$ws.Range("N45").Value = -2254
$ws.Range("H61").Value = 1901.7
$ws.Range("I61").Value = 1202.2667
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1202.2667
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -990.2666999999999
$ws.Range("N61").Value = -4424
$ws.Range("H116").Value = 27778452
$ws.Range("I116").Value = 52632016
$ws.Range("J116").Value = 941.35297
$ws.Range("K116").Value = 52632016
$ws.Range("L116").Value = 941.35297
$ws.Range("M116").Value = -52629722
$ws.Range("N116").Value = -5529.35297
$ws.Range("H132").Value = 3396.8596
$ws.Range("I132").Value = 3498.8914
$ws.Range("J132").Value = 2970.182
$ws.Range("K132").Value = 10496.6742
$ws.Range("L132").Value = 8910.545999999998
$ws.Range("M132").Value = -7966.674199999999
$ws.Range("N132").Value = -13970.546
$ws.Range("H133").Value = 35000
$ws.Range("J133").Value = 35000
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060
$ws.Range("H136").Value = 1901.7
$ws.Range("I136").Value = 1202.2667
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3606.800099999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1056.800099999999
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27778452
$ws.Range("I3").Value = 52632016
$ws.Range("J3").Value = 941.35297
$ws.Range("K3").Value = 52632016
$ws.Range("L3").Value = 941.35297
$ws.Range("M3").Value = -52631902
$ws.Range("N3").Value = -1169.35297
$ws.Range("H20").Value = 3356.7666
$ws.Range("I20").Value = 4015.6843
$ws.Range("J20").Value = 2218.6365
$ws.Range("K20").Value = 4015.6843
$ws.Range("L20").Value = 2218.6365
$ws.Range("M20").Value = -3768.6843
$ws.Range("N20").Value = -2712.6365
$ws.Range("H134").Value = 24098.295
$ws.Range("I134").Value = 32266.406
$ws.Range("J134").Value = 2316.6667
$ws.Range("K134").Value = 96799.21799999999
$ws.Range("L134").Value = 6950.000100000001
$ws.Range("M134").Value = -94264.21799999999
$ws.Range("N134").Value = -12020.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2876577.8
$ws.Range("I31").Value = 2265.303
$ws.Range("J31").Value = 6670670
$ws.Range("K31").Value = 2265.303
$ws.Range("L31").Value = 6670670
$ws.Range("M31").Value = -1970.303
$ws.Range("N31").Value = -6671260
$ws.Range("H34").Value = 2876577.8
$ws.Range("I34").Value = 2265.303
$ws.Range("J34").Value = 6670670
$ws.Range("K34").Value = 2265.303
$ws.Range("L34").Value = 6670670
$ws.Range("M34").Value = -2063.303
$ws.Range("N34").Value = -6671074
$ws.Range("H132").Value = 4168755.8
$ws.Range("I132").Value = 1452.9
$ws.Range("J132").Value = 12503361
$ws.Range("K132").Value = 4358.700000000001
$ws.Range("L132").Value = 37510083
$ws.Range("M132").Value = -1828.700000000001
$ws.Range("N132").Value = -37515143
$ws.Range("H134").Value = 1217.2858
$ws.Range("I134").Value = 1170.1666
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3510.4998
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -975.4998000000001
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 470.27274
$ws.Range("I46").Value = 200.5
$ws.Range("J46").Value = 794
$ws.Range("K46").Value = 601.5
$ws.Range("L46").Value = 2382
$ws.Range("M46").Value = -510.5
$ws.Range("N46").Value = -2564
$ws.Range("H51").Value = 625
$ws.Range("I51").Value = 450
$ws.Range("J51").Value = 800
$ws.Range("K51").Value = 1350
$ws.Range("L51").Value = 2400
$ws.Range("M51").Value = -890
$ws.Range("N51").Value = -3320
$ws.Range("H131").Value = 834.5700000000001
$ws.Range("I131").Value = 850
$ws.Range("J131").Value = 834.0928
$ws.Range("K131").Value = 2550
$ws.Range("L131").Value = 2502.2784
$ws.Range("M131").Value = 2490
$ws.Range("N131").Value = -12582.2784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3371.087
$ws.Range("I80").Value = 2349.4443
$ws.Range("J80").Value = 4027.8572
$ws.Range("K80").Value = 2349.4443
$ws.Range("L80").Value = 4027.8572
$ws.Range("M80").Value = -1351.4443
$ws.Range("N80").Value = -6023.8572
$ws.Range("H83").Value = 3371.087
$ws.Range("I83").Value = 2349.4443
$ws.Range("J83").Value = 4027.8572
$ws.Range("K83").Value = 11747.2215
$ws.Range("L83").Value = 20139.286
$ws.Range("M83").Value = -6755.2215
$ws.Range("N83").Value = -30123.286
$ws.Range("H102").Value = 1188.375
$ws.Range("I102").Value = 1078
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 1078
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 544
$ws.Range("N102").Value = -4910.6666
$ws.Range("H126").Value = 2930.6155
$ws.Range("I126").Value = 3442.5715
$ws.Range("J126").Value = 2333.3333
$ws.Range("K126").Value = 10327.7145
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("M126").Value = -7857.7145
$ws.Range("N126").Value = -11939.9999
$ws.Range("H131").Value = 21247
$ws.Range("J131").Value = 21247
$ws.Range("L131").Value = 21247
$ws.Range("N131").Value = -31327

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8462.25
$ws.Range("I132").Value = 9004.272000000001
$ws.Range("K132").Value = 27012.816
$ws.Range("M132").Value = -24482.816
$ws.Range("H133").Value = 22323.467
$ws.Range("J133").Value = 22323.467
$ws.Range("L133").Value = 22323.467
$ws.Range("N133").Value = -27383.467

